$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1902985074626866
$ws.Range("C2").Value = 0.5671641791044776
$ws.Range("J2").Value = 0.02985074626865672
$ws.Range("P2").Value = 0.1492537313432836
$ws.Range("S2").Value = 0.06343283582089553
$ws.Range("B3").Value = 0.01324503311258278
$ws.Range("C3").Value = 0.006622516556291391
$ws.Range("J3").Value = 0.03973509933774835
$ws.Range("P3").Value = 0.7814569536423841
$ws.Range("S3").Value = 0.1589403973509934
$ws.Range("J4").Value = 0.03571428571428571
$ws.Range("P4").Value = 0.5714285714285714
$ws.Range("S4").Value = 0.3928571428571428
$ws.Range("B6").Value = 0.03846153846153846
$ws.Range("D6").Value = 0.0170940170940171
$ws.Range("F6").Value = 0.07692307692307693
$ws.Range("J6").Value = 0.217948717948718
$ws.Range("O6").Value = 0.02564102564102564
$ws.Range("Q6").Value = 0.1538461538461539
$ws.Range("R6").Value = 0.07692307692307693
$ws.Range("S6").Value = 0.3931623931623932
$ws.Range("B7").Value = 0.07653061224489796
$ws.Range("D7").Value = 0.00510204081632653
$ws.Range("F7").Value = 0.08673469387755102
$ws.Range("J7").Value = 0.1479591836734694
$ws.Range("O7").Value = 0.01530612244897959
$ws.Range("Q7").Value = 0.1683673469387755
$ws.Range("R7").Value = 0.08163265306122448
$ws.Range("S7").Value = 0.4183673469387755
$ws.Range("B8").Value = 0.06553911205073996
$ws.Range("D8").Value = 0.006342494714587738
$ws.Range("E8").Value = 0.004228329809725159
$ws.Range("F8").Value = 0.05919661733615222
$ws.Range("J8").Value = 0.08879492600422834
$ws.Range("O8").Value = 0.02325581395348837
$ws.Range("Q8").Value = 0.2198731501057082
$ws.Range("R8").Value = 0.0824524312896406
$ws.Range("S8").Value = 0.4503171247357294
$ws.Range("B9").Value = 0.07216494845360824
$ws.Range("D9").Value = 0.005154639175257732
$ws.Range("F9").Value = 0.07216494845360824
$ws.Range("J9").Value = 0.09793814432989691
$ws.Range("Q9").Value = 0.211340206185567
$ws.Range("R9").Value = 0.07216494845360824
$ws.Range("S9").Value = 0.4690721649484536
$ws.Range("B10").Value = 0.1058394160583942
$ws.Range("D10").Value = 0.01532846715328467
$ws.Range("E10").Value = 0.00072992700729927
$ws.Range("F10").Value = 0.06715328467153285
$ws.Range("J10").Value = 0.1167883211678832
$ws.Range("O10").Value = 0.006569343065693431
$ws.Range("Q10").Value = 0.2503649635036496
$ws.Range("R10").Value = 0.06861313868613139
$ws.Range("S10").Value = 0.3686131386861314
$ws.Range("G11").Value = 0.1305841924398626
$ws.Range("J11").Value = 0.08247422680412371
$ws.Range("K11").Value = 0.1752577319587629
$ws.Range("L11").Value = 0.5979381443298969
$ws.Range("S11").Value = 0.01374570446735395
$ws.Range("G12").Value = 0.7634408602150538
$ws.Range("J12").Value = 0.1505376344086022
$ws.Range("K12").Value = 0.01075268817204301
$ws.Range("L12").Value = 0.04838709677419355
$ws.Range("S12").Value = 0.02688172043010753
$ws.Range("G13").Value = 0.6944444444444444
$ws.Range("J13").Value = 0.3055555555555556
$ws.Range("F15").Value = 0.02304147465437788
$ws.Range("H15").Value = 0.1290322580645161
$ws.Range("I15").Value = 0.06912442396313365
$ws.Range("J15").Value = 0.4562211981566821
$ws.Range("K15").Value = 0.06451612903225806
$ws.Range("M15").Value = 0.004608294930875576
$ws.Range("O15").Value = 0.03686635944700461
$ws.Range("S15").Value = 0.2165898617511521
$ws.Range("F16").Value = 0.01176470588235294
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.1
$ws.Range("J16").Value = 0.4235294117647059
$ws.Range("K16").Value = 0.08235294117647059
$ws.Range("O16").Value = 0.07058823529411765
$ws.Range("S16").Value = 0.1117647058823529
$ws.Range("F17").Value = 0.01811594202898551
$ws.Range("H17").Value = 0.1847826086956522
$ws.Range("I17").Value = 0.08876811594202899
$ws.Range("J17").Value = 0.4184782608695652
$ws.Range("K17").Value = 0.1050724637681159
$ws.Range("M17").Value = 0.01449275362318841
$ws.Range("N17").Value = 0.001811594202898551
$ws.Range("O17").Value = 0.05978260869565218
$ws.Range("S17").Value = 0.108695652173913
$ws.Range("F18").Value = 0.02209944751381215
$ws.Range("H18").Value = 0.1325966850828729
$ws.Range("I18").Value = 0.0718232044198895
$ws.Range("J18").Value = 0.4861878453038674
$ws.Range("K18").Value = 0.1049723756906077
$ws.Range("M18").Value = 0.005524861878453038
$ws.Range("O18").Value = 0.07734806629834254
$ws.Range("S18").Value = 0.09944751381215469
$ws.Range("F19").Value = 0.01354401805869074
$ws.Range("H19").Value = 0.2197140707298721
$ws.Range("I19").Value = 0.07825432656132431
$ws.Range("J19").Value = 0.3844996237772761
$ws.Range("K19").Value = 0.09631301730624529
$ws.Range("M19").Value = 0.02031602708803612
$ws.Range("O19").Value = 0.06847253574115876
$ws.Range("S19").Value = 0.1188863807373965
